$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column M mirrors the formatting of column L (the previous last
# column) for each populated row, then gets its own 2022 values.

$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial(-4122)

$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("M4").Value = 2022

$ws.Range("L6").Copy()
$ws.Range("M6").PasteSpecial(-4122)
$ws.Range("M6").Value = 18

$ws.Range("L7").Copy()
$ws.Range("M7").PasteSpecial(-4122)
$ws.Range("M7").Value = 6.2

$ws.Range("L8").Copy()
$ws.Range("M8").PasteSpecial(-4122)
$ws.Range("M8").Value = "-"

$ws.Range("N4").Select()
